$d = $word.ActiveDocument

# Locate the existing "{name.first} {name.last}" placeholder text (this also
# spans the spell-check proofErr-wrapped runs for "name.first"/"name.last").
$findRange = $d.Content
$found = $findRange.Find.Execute("{name.first} {name.last}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $findRange.Start

    # Remember formatting (font size, in points) so the replacement text keeps it.
    $fontSize = $findRange.Font.Size

    # Remove the old placeholder entirely - this also drops the w:proofErr
    # spell-check markers that wrapped "name.first" / "name.last", since
    # those words no longer exist in the new markup.
    $findRange.Delete()

    # Rebuild the new template markup as a sequence of inserted runs:
    #   {#name}{first}, {last}{/name}
    $segments = @("#name", "}", "{first}, {last}", "{", "/name", "}")

    $pos = $start
    $insertion = $d.Range($pos, $pos)
    $insertion.InsertAfter("{")
    $pos = $pos + 1
    $r = $d.Range($start, $pos)
    $r.Font.Size = $fontSize

    foreach ($seg in $segments) {
        $insertion = $d.Range($pos, $pos)
        $insertion.InsertAfter($seg)
        $newPos = $pos + $seg.Length
        $r = $d.Range($pos, $newPos)
        $r.Font.Size = $fontSize
        $pos = $newPos
    }
}
